$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A34").Formula = '=B17/A17'
$ws.Range("C1").Formula = '=A1*$A$34'
$ws.Range("D1").Formula = '=B1'
$ws.Range("C2:C32").Formula = '=A2*$A$34'
$ws.Range("D2:D32").Formula = '=B2'

$co2 = $ws.ChartObjects().Add(400, 400, 300, 200)
$chart2 = $co2.Chart
$chart2.ChartType = -4169
$series = $chart2.SeriesCollection()
$s1 = $series.NewSeries()
$s1.XValues = $ws.Range("C2:C32")
$s1.Values = $ws.Range("D2:D32")
$s1.MarkerStyle = 8
$chart2.HasLegend = $false
